# Insert two new data rows before the current row 440 ("Cuatro cascos verde" /
# 2021-11-05 Maule entry), shifting the existing rows 440:517 down to 442:519.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 440 (pushes old 440.. down to 442..)
$ws.Range("A440:A441").EntireRow.Insert()

# New row 440: Pimiento / Zafiro rojo / Arica y Parinacota
$ws.Range("A440").Value = 5
$ws.Range("B440").Value = "Macroferia Regional de Talca"
$ws.Range("C440").Value = "Maule"
$ws.Range("D440").Value = 44694
$ws.Range("E440").Value = 7
$ws.Range("F440").Value = 100112002
$ws.Range("G440").Value = "Pimiento"
$ws.Range("H440").Value = "Zafiro rojo"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 200
$ws.Range("K440").Value = 25000
$ws.Range("L440").Value = 25000
$ws.Range("M440").Value = 25000
$ws.Range("N440").Value = "$/caja 15 kilos"
$ws.Range("O440").Value = "Región de Arica y Parinacota"
$ws.Range("P440").Value = 1667
$ws.Range("Q440").Value = 15
$ws.Range("R440").Value = "Hortaliza"

# New row 441: Pimiento / Zafiro verde / Arica y Parinacota
$ws.Range("A441").Value = 5
$ws.Range("B441").Value = "Macroferia Regional de Talca"
$ws.Range("C441").Value = "Maule"
$ws.Range("D441").Value = 44694
$ws.Range("E441").Value = 7
$ws.Range("F441").Value = 100112002
$ws.Range("G441").Value = "Pimiento"
$ws.Range("H441").Value = "Zafiro verde"
$ws.Range("I441").Value = "Primera"
$ws.Range("J441").Value = 200
$ws.Range("K441").Value = 15000
$ws.Range("L441").Value = 15000
$ws.Range("M441").Value = 15000
$ws.Range("N441").Value = "$/caja 15 kilos"
$ws.Range("O441").Value = "Región de Arica y Parinacota"
$ws.Range("P441").Value = 1000
$ws.Range("Q441").Value = 15
$ws.Range("R441").Value = "Hortaliza"
